# Applies the "Keep.docx" diff:
#  - Reworks the "Sets:" bullet block into a new "Sets (Wrapers):" block
#    with several new bullets, reworded bullets and two bullets removed.
#  - Adds a new bullet right after "Aggregate SPO Resources from
#    Statements / Transforms / Mappings / Kinds."

$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $found = $d.Content.Find.Execute($oldText, $true, $true, $false, $false, $false, `
                                      $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found -> $oldText"
    }
}

function Find-ParagraphByText($needle) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

function Insert-BulletAfter($paragraph, $text) {
    $paragraph.Range.InsertParagraphAfter()
    $newPara = $paragraph.Next()
    $newPara.Range.Text = $text
    return $newPara
}

# --- 1. "Sets:" -> "Sets (Wrapers):" -----------------------------------
Replace-ExactText "Sets:" "Sets (Wrapers):"

# --- 2. Insert the six brand-new bullets right after it ----------------
$anchor = Find-ParagraphByText "Sets (Wrapers):"
$anchor = Insert-BulletAfter $anchor "Resource : Populate SPO Sets."
$anchor = Insert-BulletAfter $anchor "Kind: Resource Context. Aggregate Kinds."
$anchor = Insert-BulletAfter $anchor "Statement: Kind Context. Build Statements for each SPO Kind."
$anchor = Insert-BulletAfter $anchor "Mapping: Statement Context. Core Model."
$anchor = Insert-BulletAfter $anchor "Transform: Mapping Context. Core Model."
$anchor = Insert-BulletAfter $anchor "Sets (Wrapped): Dimension, Time, Measure, Employment, etc. Model reified."

# --- 3. Reword the "Universe:" bullet -----------------------------------
Replace-ExactText "Universe: Transforms: (Class / Transform, SK, PK / Mapping, OK);" `
                  "Universe: Transforms: (Mapping, SK, PK, OK);"

# --- 4. Remove the "Mappings: Kinds (in SPO Contexts: Occurrences)." bullet
$p = Find-ParagraphByText "Mappings: Kinds (in SPO Contexts: Occurrences)."
if ($p -ne $null) { $p.Range.Delete() }

# --- 5. Reword the Subjects / Predicates / Objects bullets --------------
Replace-ExactText "Subjects: (Transform / Class, Subject : Resource, Attribute : P, Value : O);" `
                  "Subjects: (SubjectKind, Subject : Resource, Attribute : P, Value : O);"

Replace-ExactText "Predicates: (Transform / Class, Attribute : S, Predicate : Resource, Value : O);" `
                  "Predicates: (PredicateKind, Attribute : S, Predicate : Resource, Value : O);"

Replace-ExactText "Objects: (Transform / Class, Attribute : P, Value : S, Object : Resource);" `
                  "Objects: (ObjectKind, Attribute : P, Value : S, Object : Resource);"

# --- 6. Reword the SubjectKind bullets -----------------------------------
Replace-ExactText "SubjectKind (SK): Predicate / Object Intersection. SubjectKind: Mapping. Occurring." `
                  "SubjectKind (SK): Predicate / Object Intersection. Occurrence:"

Replace-ExactText "(Subject : SubjectKind, Predicate : Resource, Object : Resource);" `
                  "(Context : Statement, Subject : SubjectKind, Predicate : Resource, Object : Resource);"

# --- 7. Reword the PredicateKind bullets ---------------------------------
Replace-ExactText "PredicateKind (PK / Mapping): Subject / Object intersection. PredicateKind: Mapping." `
                  "PredicateKind (PK): Subject / Object intersection:"

Replace-ExactText "(Subject : Resource, Predicate : PredicateKind, Object : Resource);" `
                  "(Context : Statement, Subject : Resource, Predicate : PredicateKind, Object : Resource);"

# --- 8. Reword the ObjectKind bullets -------------------------------------
Replace-ExactText "ObjectKind (OK): Predicate / Subject intersection. ObjectKind: Mapping." `
                  "ObjectKind (OK): Predicate / Subject intersection. Occurring."

Replace-ExactText "(Subject : Resource, Predicate : Resource, Object : ObjectKind);" `
                  "(Context : Statement, Subject : Resource, Predicate : Resource, Object : ObjectKind);"

# --- 9. Reword the Statements bullets -------------------------------------
Replace-ExactText "Statements: Subject / Predicate / Object intersection. Mapping Class." `
                  "Statements: Subject / Predicate / Object intersection:"

Replace-ExactText "(Context : Transform / Class, Subject : Resource, Predicate : Resource, Object : Resource);" `
                  "(Context : Mapping, Subject : Resource, Predicate : Resource, Object : Resource);"

# --- 10. Remove the trailing "Mappings: Context Class / Predicate / PredicateKind." bullet
$p = Find-ParagraphByText "Mappings: Context Class / Predicate / PredicateKind."
if ($p -ne $null) { $p.Range.Delete() }

# --- 11. Add the new bullet after "Aggregate SPO Resources ..." ----------
$p = Find-ParagraphByText "Aggregate SPO Resources from Statements / Transforms / Mappings / Kinds."
if ($p -ne $null) {
    Insert-BulletAfter $p "Align Statements with Core Model: Templates / Transforms."
}

Write-Host "Done"
